$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrease trial_total (column F) by 78 for each data row (rows 2-42)
# so that n distractors = n targets
for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value2 = $cell.Value2 - 78
}
